# Apply updated dSF (column F) values as per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -5
    "F6"  = -1
    "F7"  = -1
    "F8"  = -1
    "F14" = -4
    "F17" = -7
    "F19" = -9
    "F20" = 3
    "F23" = 8
    "F26" = 4
    "F28" = -9
    "F29" = 1
    "F31" = -4
    "F32" = -2
    "F33" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
